$d = $word.ActiveDocument

# 1. Append " (avec mode silence)" after "Eisenhower Matrix." in the productivity methods bullet.
$d.Content.Find.Execute(", Eisenhower Matrix.", $true, $false, $false, $false, $false, $true, 1, $false, ", Eisenhower Matrix. (avec mode silence)", 2)

# 2. Collapse "Rida Elantari." (previously split across runs with spellcheck markers) into a single run's text.
$d.Content.Find.Execute("Rida Elantari.", $true, $false, $false, $false, $false, $true, 1, $false, "Rida Elantari.", 2)
